# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Cant. Trabajadores" (count of workers) table used to list 3 workers
# across several "Periodo Mora" rows (7 data rows, rows 16-22). It is
# refreshed to list 2 workers (JEINER GARCIA VIDES / ORLANY MANUEL VERGARA
# BELTRAN) each with 2 periods (2507, 2508) -> 4 data rows (16-19), and the
# summary figures (Valor Mora, Cant. Trabajadores, Cant. Periodos) are
# updated to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Copy the "bottom border" formatting of the old last data row (22) onto
#    what will become the new last data row (19) before we touch anything,
#    so the table keeps its closing border once the extra rows are removed.
$ws.Range("B22:J22").Copy()
$ws.Range("B19:J19").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# 2) Refresh the worker / period detail table (rows 16-19).
#    Row 16: JEINER GARCIA VIDES, periodo 2507
$ws.Range("C16").Value = "1051676899"
$ws.Range("D16").Value = "JEINER GARCIA VIDES"
$ws.Range("E16").Value = "2507"

#    Row 17: ORLANY MANUEL VERGARA BELTRAN, periodo 2507
$ws.Range("C17").Value = "1052524282"
$ws.Range("D17").Value = "ORLANY MANUEL VERGARA BELTRAN"
$ws.Range("E17").Value = "2507"

#    Row 18: JEINER GARCIA VIDES, periodo 2508
$ws.Range("C18").Value = "1051676899"
$ws.Range("D18").Value = "JEINER GARCIA VIDES"
$ws.Range("E18").Value = "2508"

#    Row 19: ORLANY MANUEL VERGARA BELTRAN, periodo 2508
$ws.Range("C19").Value = "1052524282"
$ws.Range("D19").Value = "ORLANY MANUEL VERGARA BELTRAN"
$ws.Range("E19").Value = "2508"
$ws.Range("F19").Value = 56940

# 3) Remove the now-unused trailing data rows (old rows 20-22). This shifts
#    the signature block below (old rows 27-28) up to rows 24-25.
$ws.Rows("20:22").Delete()

# 4) Update the summary figures to match the refreshed table.
$ws.Range("E11").Value = 227760   # VALOR MORA
$ws.Range("C13").Value = 2        # Cant. Trabajadores
$ws.Range("F13").Value = 2        # Cant. Periodos
